$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.931.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.644.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.58%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5097"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.56%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2575"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06413"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07756"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.72%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.308"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.55%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.658.96"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.47%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5467"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅7887"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.41%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.995.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.50%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.005"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "198.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.439"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.67%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.056"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.31%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.006"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.852"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1146"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.29%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.895"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.88%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.44%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.239"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05016"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.62%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.284"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.83%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.202"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.539"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.364"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8933"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.580"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.132.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5532"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01563"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.005"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.658"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8149"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.87%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.52%  "

$ws.Range("E44").Value = "  +8.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.785.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4525"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.20%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.95%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.34%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05091"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.09568"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.35%  "

$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.006"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.11%  "
